$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Periodo Mora" / "Valor Mora" pair between row 16 and row 18.
$ws.Range("E16").Value = "2309"
$ws.Range("F16").Value = 46400
$ws.Range("E18").Value = "2406"
$ws.Range("F18").Value = 20800
